# Updated cryptos list on Tue Dec  5 16:32:26 UTC 2023 with GitHub Actions
#
# Applies the refreshed price / 1h-volume figures (and the re-ranking of a
# few coins that swapped positions) to Sheet1 of the cryptos workbook.
#
# All of the touched cells are stored as plain text in the source workbook
# (prices/percentages are formatted strings, not numbers), so every write
# goes through Set-TextCell, which forces the cell to Text format before
# assigning the value (otherwise numeric-looking strings like "231.44"
# would silently become real numbers) and then restores the cell's style
# to Normal so no stray formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Bitcoin
Set-TextCell 2 4 "42.466.69"
Set-TextCell 2 5 "  +2.35%  "
# Ethereum
Set-TextCell 3 4 "2.227.87"
Set-TextCell 3 5 "  +0.31%  "
# TetherUSD
Set-TextCell 4 5 "  -0.03%  "
# BNB
Set-TextCell 5 4 "231.44"
Set-TextCell 5 5 "  +0.09%  "
# XRP
Set-TextCell 6 4 "0.621"
Set-TextCell 6 5 "  -0.55%  "
# Solana
Set-TextCell 7 4 "60.96"
Set-TextCell 7 5 "  -0.22%  "
# USDC
Set-TextCell 8 5 "  -0.01%  "
# Cardano
Set-TextCell 9 4 "0.405"
Set-TextCell 9 5 "  +0.70%  "
# Dogecoin
Set-TextCell 10 4 "0.0909"
Set-TextCell 10 5 "  +2.08%  "
# TRON
Set-TextCell 11 5 "  +0.02%  "
# WrappedliquidstakedEther2.0
Set-TextCell 12 4 "2.559.45"
Set-TextCell 12 5 "  +0.56%  "
# Chainlink
Set-TextCell 13 5 "  -0.75%  "
# Avalanche
Set-TextCell 14 4 "22.34"
Set-TextCell 14 5 "  +2.70%  "
# Polkadot
Set-TextCell 15 4 "5.65"
Set-TextCell 15 5 "  +1.86%  "
# Polygon
Set-TextCell 16 5 "  +0.09%  "
# WrappedEther
Set-TextCell 17 4 "2.245.98"
Set-TextCell 17 5 "  +1.30%  "
# WrappedBTC
Set-TextCell 18 4 "42.285.13"
Set-TextCell 18 5 "  +2.33%  "
# ShibaInu
Set-TextCell 19 4 "0.0₃0948"
Set-TextCell 19 5 "  +5.78%  "
# row 20
Set-TextCell 20 5 "  +2.06%  "
# row 21
Set-TextCell 21 4 "72.29"
Set-TextCell 21 5 "  -0.80%  "
# row 22
Set-TextCell 22 4 "244.49"
Set-TextCell 22 5 "  -2.10%  "
# row 23
Set-TextCell 23 5 "  -0.16%  "
# row 24
Set-TextCell 24 4 "2.46"
Set-TextCell 24 5 "  +3.07%  "
# row 25
Set-TextCell 25 4 "2.31"
Set-TextCell 25 5 "  +1.56%  "
# row 26
Set-TextCell 26 4 "9.71"
Set-TextCell 26 5 "  +2.90%  "
# row 27
Set-TextCell 27 4 "169.29"
Set-TextCell 27 5 "  +0.60%  "
# row 28
Set-TextCell 28 4 "0.141"
Set-TextCell 28 5 "  +1.89%  "
# row 29
Set-TextCell 29 4 "20.42"
Set-TextCell 29 5 "  +2.35%  "
# row 30
Set-TextCell 30 5 "  +2.54%  "
# row 31
Set-TextCell 31 4 "2.65"
Set-TextCell 31 5 "  +0.84%  "
# row 32
Set-TextCell 32 5 "  -1.66%  "

# Row 33 / 34 swap: Filecoin now ranks above InternetComputer(DFINITY)
Set-TextCell 33 2 "Filecoin"
Set-TextCell 33 3 "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell 33 4 "4.79"
Set-TextCell 33 5 "  +3.33%  "
Set-TextCell 34 2 "InternetComputer(DFINITY)"
Set-TextCell 34 3 "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell 34 4 "5.00"
Set-TextCell 34 5 "  -0.12%  "

# row 35
Set-TextCell 35 4 "0.0654"
Set-TextCell 35 5 "  +4.74%  "
# row 36
Set-TextCell 36 4 "6.39"
Set-TextCell 36 5 "  -2.39%  "
# row 37
Set-TextCell 37 4 "2.37"
Set-TextCell 37 5 "  +0.06%  "
# row 38
Set-TextCell 38 4 "3.58"
Set-TextCell 38 5 "  -2.98%  "
# row 39
Set-TextCell 39 4 "0.0249"
Set-TextCell 39 5 "  +4.83%  "
# row 40
Set-TextCell 40 5 "  -0.05%  "
# row 41
Set-TextCell 41 4 "8.68"
Set-TextCell 41 5 "  +1.18%  "
# row 42
Set-TextCell 42 4 "0.000227"
Set-TextCell 42 5 "  -8.60%  "
# row 43
Set-TextCell 43 4 "0.0965"
Set-TextCell 43 5 "  -1.78%  "
# row 44
Set-TextCell 44 5 "  +0.52%  "
# row 45
Set-TextCell 45 4 "97.34"
Set-TextCell 45 5 "  -1.82%  "
# row 46
Set-TextCell 46 4 "1.457.20"
Set-TextCell 46 5 "  -0.53%  "
# row 47
Set-TextCell 47 4 "4.37"
Set-TextCell 47 5 "  -9.53%  "

# Rows 48 / 49 / 50 rotate: InjectiveProtocol, HuobiToken, ARBITRUM shift up one slot
Set-TextCell 48 2 "InjectiveProtocol"
Set-TextCell 48 3 "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell 48 4 "16.17"
Set-TextCell 48 5 "  -2.67%  "
Set-TextCell 49 2 "HuobiToken"
Set-TextCell 49 3 "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell 49 4 "2.75"
Set-TextCell 49 5 "  -1.05%  "
Set-TextCell 50 2 "ARBITRUM"
Set-TextCell 50 3 "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell 50 4 "1.07"
Set-TextCell 50 5 "  -0.55%  "

# row 51
Set-TextCell 51 4 "2.21"
Set-TextCell 51 5 "  +3.85%  "
